$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$ws3 = $wb.Worksheets.Item("AVLo-passengers")
$ws4 = $wb.Worksheets.Item("AVLo-freight")

# Insert a new row after existing row 36 (pushes row 37+ down by one)
$ws2.Rows.Item(37).Insert() | Out-Null

# New row 37: weighted value note + formula referencing the AVg loading value above
$ws2.Range("A37").Value = "weighted value, adjusted for number of train cars per locomotive"
$ws2.Range("A37").WrapText = $true

$ws2.Range("B37").Formula = "=B36/10"
$ws2.Range("B37").NumberFormat = "0"

# Row heights for the (now taller) summary row and the new row
$ws2.Rows.Item(36).RowHeight = 16
$ws2.Rows.Item(37).RowHeight = 16

# Fix up references on the other sheets that pointed at rows which shifted down
$ws3.Range("B5").Formula = "='BTS NTS Modal Profile Data'!B37"
$ws3.Range("B7").Formula = "='BTS NTS Modal Profile Data'!B60"
$ws4.Range("B6").Formula = "='BTS NTS Modal Profile Data'!B55"

# Restore/update selections so the saved workbook view matches
$ws1.Range("A44").Select() | Out-Null

$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 14
$ws2.Range("C34").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("F12").Select() | Out-Null
